# Add data for 2024-12-29
# Updates the 2024 YTD (column K) violent-crime figures across the
# Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item(1)
$ws.Range("K2").Value = 7849
$ws.Range("K3").Value = 8136
$ws.Range("I4").Value = 1809
$ws.Range("K4").Value = 1714
$ws.Range("K6").Value = 9055
$ws.Range("I7").Value = 26057
$ws.Range("K7").Value = 27332

# Sheet 11: Logan Square
$ws = $wb.Worksheets.Item(11)
$ws.Range("K6").Value = 145
$ws.Range("K7").Value = 347

# Sheet 12: Austin
$ws = $wb.Worksheets.Item(12)
$ws.Range("K3").Value = 542
$ws.Range("K6").Value = 602
$ws.Range("K7").Value = 1788

# Sheet 13: South Chicago
$ws = $wb.Worksheets.Item(13)
$ws.Range("K2").Value = 202
$ws.Range("K7").Value = 582

# Sheet 15: West Pullman
$ws = $wb.Worksheets.Item(15)
$ws.Range("K2").Value = 153
$ws.Range("K7").Value = 454

# Sheet 16: Grand Crossing
$ws = $wb.Worksheets.Item(16)
$ws.Range("K2").Value = 258
$ws.Range("K6").Value = 270
$ws.Range("K7").Value = 900

# Sheet 17: New City
$ws = $wb.Worksheets.Item(17)
$ws.Range("K2").Value = 206
$ws.Range("K3").Value = 150
$ws.Range("K7").Value = 634

# Sheet 18: Woodlawn
$ws = $wb.Worksheets.Item(18)
$ws.Range("K2").Value = 124
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 464

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item(2)
$ws.Range("K2").Value = 237
$ws.Range("K7").Value = 812
$ws.Range("K8").Value = 1788
$ws.Range("K9").Value = 131
$ws.Range("K11").Value = 479
$ws.Range("K16").Value = 66
$ws.Range("K19").Value = 784
$ws.Range("K22").Value = 85
$ws.Range("K23").Value = 267
$ws.Range("K29").Value = 1506
$ws.Range("K34").Value = 155
$ws.Range("K37").Value = 900
$ws.Range("K41").Value = 180
$ws.Range("K42").Value = 1016
$ws.Range("K48").Value = 347
$ws.Range("K50").Value = 122
$ws.Range("K51").Value = 354
$ws.Range("K53").Value = 347
$ws.Range("K54").Value = 533
$ws.Range("K55").Value = 298
$ws.Range("I63").Value = 240
$ws.Range("K63").Value = 80
$ws.Range("K65").Value = 634
$ws.Range("K67").Value = 1065
$ws.Range("K70").Value = 50
$ws.Range("K75").Value = 90
$ws.Range("K76").Value = 376
$ws.Range("K77").Value = 180
$ws.Range("K80").Value = 104
$ws.Range("K83").Value = 582
$ws.Range("K84").Value = 223
$ws.Range("K85").Value = 1266
$ws.Range("K90").Value = 260
$ws.Range("K93").Value = 111
$ws.Range("K94").Value = 365
$ws.Range("K95").Value = 454
$ws.Range("K96").Value = 296
$ws.Range("K97").Value = 225
$ws.Range("K98").Value = 147
$ws.Range("K99").Value = 464
$ws.Range("I101").Value = 26057
$ws.Range("K101").Value = 27332

# Sheet 21: North Lawndale
$ws = $wb.Worksheets.Item(21)
$ws.Range("K2").Value = 292
$ws.Range("K3").Value = 387
$ws.Range("K4").Value = 61
$ws.Range("K7").Value = 1065

# Sheet 22: South Deering
$ws = $wb.Worksheets.Item(22)
$ws.Range("K3").Value = 90
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 223

# Sheet 24: Loop
$ws = $wb.Worksheets.Item(24)
$ws.Range("K4").Value = 39
$ws.Range("K6").Value = 287
$ws.Range("K7").Value = 533

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item(25)
$ws.Range("K3").Value = 532
$ws.Range("K7").Value = 1506

# Sheet 26: Lake View
$ws = $wb.Worksheets.Item(26)
$ws.Range("K3").Value = 83
$ws.Range("K4").Value = 52
$ws.Range("K6").Value = 158
$ws.Range("K7").Value = 347

# Sheet 27: Chatham
$ws = $wb.Worksheets.Item(27)
$ws.Range("K6").Value = 262
$ws.Range("K7").Value = 784

# Sheet 29: River North
$ws = $wb.Worksheets.Item(29)
$ws.Range("K6").Value = 185
$ws.Range("K7").Value = 376

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item(31)
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 180

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item(32)
$ws.Range("K2").Value = 270
$ws.Range("K3").Value = 297
$ws.Range("K6").Value = 388
$ws.Range("K7").Value = 1016

# Sheet 36: Lower West Side
$ws = $wb.Worksheets.Item(36)
$ws.Range("K6").Value = 111
$ws.Range("K7").Value = 298

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item(39)
$ws.Range("K3").Value = 93
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 267

# Sheet 4: West Ridge
$ws = $wb.Worksheets.Item(4)
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 296

# Sheet 48: West Lawn
$ws = $wb.Worksheets.Item(48)
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 111

# Sheet 5: Auburn Gresham
$ws = $wb.Worksheets.Item(5)
$ws.Range("K3").Value = 257
$ws.Range("K7").Value = 812

# Sheet 50: Garfield Ridge
$ws = $wb.Worksheets.Item(50)
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 155

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item(51)
$ws.Range("K3").Value = 75
$ws.Range("K6").Value = 171
$ws.Range("K7").Value = 365

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item(55)
$ws.Range("K2").Value = 28
$ws.Range("K7").Value = 147

# Sheet 56: Lincoln Square
$ws = $wb.Worksheets.Item(56)
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 122

# Sheet 6: Belmont Cragin
$ws = $wb.Worksheets.Item(6)
$ws.Range("K4").Value = 28
$ws.Range("K7").Value = 479

# Sheet 61: Avalon Park
$ws = $wb.Worksheets.Item(61)
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 131

# Sheet 64: Albany Park
$ws = $wb.Worksheets.Item(64)
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 237

# Sheet 65: West Town
$ws = $wb.Worksheets.Item(65)
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 126
$ws.Range("K7").Value = 225

# Sheet 67: O'Hare
$ws = $wb.Worksheets.Item(67)
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 50

# Sheet 73: Pullman
$ws = $wb.Worksheets.Item(73)
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 90

# Sheet 74: Washington Heights
$ws = $wb.Worksheets.Item(74)
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 260

# Sheet 75: Little Italy, UIC
$ws = $wb.Worksheets.Item(75)
$ws.Range("K3").Value = 98
$ws.Range("K6").Value = 115
$ws.Range("K7").Value = 354

# Sheet 8: South Shore
$ws = $wb.Worksheets.Item(8)
$ws.Range("K2").Value = 421
$ws.Range("K3").Value = 437
$ws.Range("K6").Value = 311
$ws.Range("K7").Value = 1266

# Sheet 80: Clearing
$ws = $wb.Worksheets.Item(80)
$ws.Range("K6").Value = 19
$ws.Range("K7").Value = 85

# Sheet 84: Riverdale
$ws = $wb.Worksheets.Item(84)
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 180

# Sheet 87: Rush & Division
$ws = $wb.Worksheets.Item(87)
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 104

# Sheet 94: Bucktown
$ws = $wb.Worksheets.Item(94)
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 66
